$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Step 1: increment the date in A1 by one day
$ws.Range("A1").Value = 45309

# Step 2: update price column D for rows 35-39 (COMUN table) and 42-46 (CON TOPE table)
$ws.Range("D35").Value = 3030.119
$ws.Range("D36").Value = 3310
$ws.Range("D37").Value = 3310
$ws.Range("D38").Value = 3310
$ws.Range("D39").Value = 8628.17

$ws.Range("D42").Value = 3267.016
$ws.Range("D43").Value = 3500
$ws.Range("D44").Value = 3500
$ws.Range("D45").Value = 3500
$ws.Range("D46").Value = 9649.955
